# Update "想去人数" (interest count) figures in the source sheets.
# Corresponds to a data refresh (gh-pages build at commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5882
$ws1.Range("F10").Value = 17
$ws1.Range("F15").Value = 1539
$ws1.Range("F19").Value = 4402
$ws1.Range("F34").Value = 1200

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 638
$ws3.Range("F5").Value = 262

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 638
$ws4.Range("F8").Value = 5882
$ws4.Range("F10").Value = 73
$ws4.Range("F22").Value = 17
$ws4.Range("F26").Value = 1539
$ws4.Range("F30").Value = 4402
